$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 00:20"

# Row 7 - China: new cases / deaths-today reset to 0
$ws.Range("C7").Value = 0
$ws.Range("G7").Value = 0

# Row 8 - Alemania: updated totals
$ws.Range("B8").Value = 71808
$ws.Range("C8").Value = 4923
$ws.Range("E8").Value = 54933
$ws.Range("G8").Value = 130
$ws.Range("H8").Value = 775

# Rows 22/23 - Australia overtakes Noruega in the ranking (rows swap + data update)
$ws.Range("A22").Value = "Australia"
$ws.Range("B22").Value = 4711
$ws.Range("C22").Value = 251
$ws.Range("D22").Value = 337
$ws.Range("E22").Value = 4354
$ws.Range("F22").Value = 28
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 20

$ws.Range("A23").Value = "Noruega"
$ws.Range("B23").Value = 4641
$ws.Range("C23").Value = 196
$ws.Range("D23").Value = 13
$ws.Range("E23").Value = 4589
$ws.Range("F23").Value = 97
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 39
